$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -7.865
$ws.Range("D4").Value = -7.769
$ws.Range("C11").Value = -12.379
$ws.Range("C12").Value = -10.896
$ws.Range("D14").Value = -7.722
$ws.Range("C15").Value = -13.117
$ws.Range("D26").Value = -8.026999999999999
$ws.Range("C27").Value = -13.046
$ws.Range("C28").Value = -13.012
$ws.Range("C31").Value = -13.647
$ws.Range("D31").Value = -8.51
$ws.Range("C32").Value = -13.895
$ws.Range("D35").Value = -7.854000000000001
$ws.Range("C36").Value = -12.732
$ws.Range("D37").Value = -7.632
$ws.Range("C38").Value = -12.703
$ws.Range("D39").Value = -7.217999999999999
$ws.Range("D40").Value = -7.854000000000001
$ws.Range("D45").Value = -7.772999999999999
$ws.Range("C46").Value = -14.009
$ws.Range("D52").Value = -7.247999999999999
$ws.Range("C54").Value = -13.408
$ws.Range("C55").Value = -13.307
$ws.Range("C56").Value = -13.364
$ws.Range("D57").Value = -8.451000000000001
$ws.Range("C67").Value = -11.705
$ws.Range("C69").Value = -11.038
$ws.Range("C72").Value = -11.555
$ws.Range("C73").Value = -13.013
$ws.Range("D81").Value = -7.221000000000001
$ws.Range("C83").Value = -13.314
$ws.Range("D83").Value = -8.461
$ws.Range("C86").Value = -13.846
$ws.Range("C91").Value = -10.885
$ws.Range("C93").Value = -11.979
$ws.Range("C99").Value = -12.635
$ws.Range("D100").Value = -8.200999999999999
$ws.Range("D102").Value = -7.577000000000001
